# Auto-generated script applying scheduled market-data refresh to Phoenix_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across all 8 sheets
$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 713.63635
$ws.Range("I5").Value = 651.6667
$ws.Range("J5").Value = 992.5
$ws.Range("K5").Value = 651.6667
$ws.Range("L5").Value = 992.5
$ws.Range("M5").Value = -536.6667
$ws.Range("N5").Value = -1222.5
$ws.Range("H8").Value = 397
$ws.Range("I8").Value = 397
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1191
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1052
$ws.Range("N8").ClearContents()
$ws.Range("H12").Value = 3875
$ws.Range("H15").Value = 680.3889
$ws.Range("I15").Value = 680.3889
$ws.Range("K15").Value = 2041.1667
$ws.Range("M15").Value = -1872.1667
$ws.Range("H61").Value = 115
$ws.Range("I61").Value = 115
$ws.Range("K61").Value = 345
$ws.Range("M61").Value = -173
$ws.Range("H62").Value = 3265.5
$ws.Range("I62").Value = 3196.125
$ws.Range("J62").Value = 3543
$ws.Range("K62").Value = 3196.125
$ws.Range("L62").Value = 3543
$ws.Range("M62").Value = -2572.125
$ws.Range("N62").Value = -4791
$ws.Range("H65").Value = 3265.5
$ws.Range("I65").Value = 3196.125
$ws.Range("J65").Value = 3543
$ws.Range("K65").Value = 15980.625
$ws.Range("L65").Value = 17715
$ws.Range("M65").Value = -12860.625
$ws.Range("N65").Value = -23955
$ws.Range("H76").Value = 8140.2104
$ws.Range("I76").Value = 8657.308000000001
$ws.Range("K76").Value = 8657.308000000001
$ws.Range("M76").Value = -8342.308000000001
$ws.Range("H79").Value = 8140.2104
$ws.Range("I79").Value = 8657.308000000001
$ws.Range("K79").Value = 8657.308000000001
$ws.Range("M79").Value = -7565.308000000001
$ws.Range("H98").Value = 61129.562
$ws.Range("I98").Value = 74051.30499999999
$ws.Range("K98").Value = 74051.30499999999
$ws.Range("M98").Value = -72553.30499999999
$ws.Range("H100").Value = 3049.2354
$ws.Range("I100").Value = 2995.2307
$ws.Range("J100").Value = 3224.75
$ws.Range("K100").Value = 2995.2307
$ws.Range("L100").Value = 3224.75
$ws.Range("M100").Value = -2454.2307
$ws.Range("N100").Value = -4306.75
$ws.Range("H106").Value = 2322.0881
$ws.Range("I106").Value = 2335.7666
$ws.Range("K106").Value = 2335.7666
$ws.Range("M106").Value = -1704.7666
$ws.Range("H112").Value = 1545.52
$ws.Range("J112").Value = 1898.7858
$ws.Range("L112").Value = 5696.357400000001
$ws.Range("N112").Value = -7912.357400000001
$ws.Range("H113").Value = 3333
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -9008
$ws.Range("H116").Value = 7554.9443
$ws.Range("I116").Value = 7299.16
$ws.Range("J116").Value = 8136.273
$ws.Range("K116").Value = 7299.16
$ws.Range("L116").Value = 8136.273
$ws.Range("M116").Value = -3857.16
$ws.Range("N116").Value = -15020.273
$ws.Range("H122").Value = 61129.562
$ws.Range("I122").Value = 74051.30499999999
$ws.Range("K122").Value = 222153.915
$ws.Range("M122").Value = -219703.915
$ws.Range("H132").Value = 3962.375
$ws.Range("I132").Value = 3962.375
$ws.Range("K132").Value = 11887.125
$ws.Range("M132").Value = -9357.125
$ws.Range("H135").Value = 1432.0869
$ws.Range("I135").Value = 1432.0869
$ws.Range("K135").Value = 12888.7821
$ws.Range("M135").Value = -10353.7821
$ws.Range("H137").Value = 144124.64
$ws.Range("I137").Value = 1435.9615
$ws.Range("J137").Value = 320786.8
$ws.Range("K137").Value = 4307.8845
$ws.Range("L137").Value = 962360.3999999999
$ws.Range("M137").Value = -1757.8845
$ws.Range("N137").Value = -967460.3999999999
$ws.Range("H138").Value = 1958.8431
$ws.Range("I138").Value = 1031.862
$ws.Range("J138").Value = 3180.7727
$ws.Range("K138").Value = 3095.586
$ws.Range("L138").Value = 9542.3181
$ws.Range("M138").Value = 2044.414
$ws.Range("N138").Value = -19822.3181
$ws.Range("H141").Value = 7251.476
$ws.Range("I141").Value = 6830
$ws.Range("J141").Value = 9042.75
$ws.Range("K141").Value = 20490
$ws.Range("L141").Value = 27128.25
$ws.Range("M141").Value = -15310
$ws.Range("N141").Value = -37488.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2092.6287
$ws.Range("I2").Value = 2068.8147
$ws.Range("J2").Value = 2173
$ws.Range("K2").Value = 2068.8147
$ws.Range("L2").Value = 2173
$ws.Range("M2").Value = -1955.8147
$ws.Range("N2").Value = -2399
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H32").Value = 3927.4614
$ws.Range("I32").Value = 3374.239
$ws.Range("K32").Value = 3374.239
$ws.Range("M32").Value = -3087.239
$ws.Range("H44").Value = 39025
$ws.Range("J44").Value = 39025
$ws.Range("L44").Value = 39025
$ws.Range("N44").Value = -40001
$ws.Range("H45").Value = 3306.6667
$ws.Range("I45").Value = 3274.2354
$ws.Range("K45").Value = 3274.2354
$ws.Range("M45").Value = -2897.2354
$ws.Range("H61").Value = 3658.6128
$ws.Range("I61").Value = 2980
$ws.Range("K61").Value = 2980
$ws.Range("M61").Value = -2768
$ws.Range("H74").Value = 59770.547
$ws.Range("I74").Value = 46408.05
$ws.Range("J74").Value = 84066
$ws.Range("K74").Value = 46408.05
$ws.Range("L74").Value = 84066
$ws.Range("M74").Value = -45534.05
$ws.Range("N74").Value = -85814
$ws.Range("H77").Value = 59770.547
$ws.Range("I77").Value = 46408.05
$ws.Range("J77").Value = 84066
$ws.Range("K77").Value = 232040.25
$ws.Range("L77").Value = 420330
$ws.Range("M77").Value = -227672.25
$ws.Range("N77").Value = -429066
$ws.Range("H80").Value = 54985.5
$ws.Range("J80").Value = 54985.5
$ws.Range("L80").Value = 54985.5
$ws.Range("N80").Value = -56981.5
$ws.Range("H83").Value = 54985.5
$ws.Range("J83").Value = 54985.5
$ws.Range("L83").Value = 164956.5
$ws.Range("N83").Value = -174940.5
$ws.Range("H97").Value = 1431.5555
$ws.Range("I97").Value = 340.7143
$ws.Range("J97").Value = 5249.5
$ws.Range("K97").Value = 340.7143
$ws.Range("L97").Value = 5249.5
$ws.Range("M97").Value = 155.2857
$ws.Range("N97").Value = -6241.5
$ws.Range("H116").Value = 2092.6287
$ws.Range("I116").Value = 2068.8147
$ws.Range("J116").Value = 2173
$ws.Range("K116").Value = 2068.8147
$ws.Range("L116").Value = 2173
$ws.Range("M116").Value = 225.1853000000001
$ws.Range("N116").Value = -6761
$ws.Range("H122").Value = 41169.957
$ws.Range("I122").Value = 1611.1333
$ws.Range("K122").Value = 4833.3999
$ws.Range("M122").Value = -2383.3999
$ws.Range("H124").Value = 41536.715
$ws.Range("J124").Value = 41536.715
$ws.Range("L124").Value = 41536.715
$ws.Range("N124").Value = -51356.715
$ws.Range("H132").Value = 8984.146000000001
$ws.Range("I132").Value = 9022.368
$ws.Range("K132").Value = 27067.104
$ws.Range("M132").Value = -24537.104
$ws.Range("H136").Value = 3658.6128
$ws.Range("I136").Value = 2980
$ws.Range("K136").Value = 8940
$ws.Range("M136").Value = -6390
$ws.Range("H137").Value = 116909.14
$ws.Range("J137").Value = 116909.14
$ws.Range("L137").Value = 116909.14
$ws.Range("N137").Value = -127109.14

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2092.6287
$ws.Range("I3").Value = 2068.8147
$ws.Range("J3").Value = 2173
$ws.Range("K3").Value = 2068.8147
$ws.Range("L3").Value = 2173
$ws.Range("M3").Value = -1954.8147
$ws.Range("N3").Value = -2401
$ws.Range("H20").Value = 8373.4
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 8373.4
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 8373.4
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -8867.4
$ws.Range("H86").Value = 2609.2415
$ws.Range("I86").Value = 2692.8235
$ws.Range("J86").Value = 2490.8333
$ws.Range("K86").Value = 2692.8235
$ws.Range("L86").Value = 2490.8333
$ws.Range("M86").Value = -1569.8235
$ws.Range("N86").Value = -4736.8333
$ws.Range("H89").Value = 2609.2415
$ws.Range("I89").Value = 2692.8235
$ws.Range("J89").Value = 2490.8333
$ws.Range("K89").Value = 13464.1175
$ws.Range("L89").Value = 12454.1665
$ws.Range("M89").Value = -7848.1175
$ws.Range("N89").Value = -23686.1665
$ws.Range("H99").Value = 2312.2
$ws.Range("J99").Value = 1876.6666
$ws.Range("L99").Value = 1876.6666
$ws.Range("N99").Value = -4872.6666
$ws.Range("H105").Value = 62501572
$ws.Range("I105").Value = 62501572
$ws.Range("K105").Value = 62501572
$ws.Range("M105").Value = -62499825
$ws.Range("H107").Value = 8621.666999999999
$ws.Range("I107").Value = 11990.125
$ws.Range("J107").Value = 1884.75
$ws.Range("K107").Value = 11990.125
$ws.Range("L107").Value = 1884.75
$ws.Range("M107").Value = -10070.125
$ws.Range("N107").Value = -5724.75
$ws.Range("H134").Value = 3492.7727
$ws.Range("I134").Value = 3676.75
$ws.Range("J134").Value = 3002.1667
$ws.Range("K134").Value = 11030.25
$ws.Range("L134").Value = 9006.500100000001
$ws.Range("M134").Value = -8495.25
$ws.Range("N134").Value = -14076.5001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 259.38095
$ws.Range("J7").Value = 271.75
$ws.Range("L7").Value = 271.75
$ws.Range("N7").Value = -497.75
$ws.Range("H15").Value = 5153
$ws.Range("J15").Value = 5153
$ws.Range("L15").Value = 5153
$ws.Range("N15").Value = -5493
$ws.Range("H22").Value = 1714
$ws.Range("I22").Value = 1499.5
$ws.Range("K22").Value = 1499.5
$ws.Range("M22").Value = -1149.5
$ws.Range("H31").Value = 1635
$ws.Range("I31").Value = 1650.4
$ws.Range("J31").Value = 1558
$ws.Range("K31").Value = 1650.4
$ws.Range("L31").Value = 1558
$ws.Range("M31").Value = -1355.4
$ws.Range("N31").Value = -2148
$ws.Range("H34").Value = 1635
$ws.Range("I34").Value = 1650.4
$ws.Range("J34").Value = 1558
$ws.Range("K34").Value = 1650.4
$ws.Range("L34").Value = 1558
$ws.Range("M34").Value = -1448.4
$ws.Range("N34").Value = -1962
$ws.Range("H41").Value = 37249.5
$ws.Range("I41").Value = 13999
$ws.Range("K41").Value = 13999
$ws.Range("M41").Value = -13571
$ws.Range("H47").Value = 35000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H50").Value = 59980.8
$ws.Range("J50").Value = 59979.5
$ws.Range("L50").Value = 59979.5
$ws.Range("N50").Value = -61229.5
$ws.Range("H51").Value = 44348
$ws.Range("J51").Value = 44348
$ws.Range("L51").Value = 44348
$ws.Range("N51").Value = -45820
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 9100.157999999999
$ws.Range("I58").Value = 3036.7273
$ws.Range("J58").Value = 17437.375
$ws.Range("K58").Value = 3036.7273
$ws.Range("L58").Value = 17437.375
$ws.Range("M58").Value = -2833.7273
$ws.Range("N58").Value = -17843.375
$ws.Range("H60").Value = 41078.582
$ws.Range("I60").Value = 14331.667
$ws.Range("K60").Value = 14331.667
$ws.Range("M60").Value = -13820.667
$ws.Range("H61").Value = 44348
$ws.Range("J61").Value = 44348
$ws.Range("L61").Value = 44348
$ws.Range("N61").Value = -45044
$ws.Range("H62").Value = 175448
$ws.Range("I62").Value = 339896.66
$ws.Range("J62").Value = 10999.333
$ws.Range("K62").Value = 339896.66
$ws.Range("L62").Value = 10999.333
$ws.Range("M62").Value = -339272.66
$ws.Range("N62").Value = -12247.333
$ws.Range("H65").Value = 175448
$ws.Range("I65").Value = 339896.66
$ws.Range("J65").Value = 10999.333
$ws.Range("K65").Value = 1699483.3
$ws.Range("L65").Value = 54996.665
$ws.Range("M65").Value = -1696363.3
$ws.Range("N65").Value = -61236.665
$ws.Range("H86").Value = 14746.5
$ws.Range("J86").Value = 16328.667
$ws.Range("L86").Value = 16328.667
$ws.Range("N86").Value = -18574.667
$ws.Range("H89").Value = 14746.5
$ws.Range("J89").Value = 16328.667
$ws.Range("L89").Value = 81643.33499999999
$ws.Range("N89").Value = -92875.33499999999
$ws.Range("H94").Value = 1340.28
$ws.Range("I94").Value = 1526
$ws.Range("J94").Value = 1168.8462
$ws.Range("K94").Value = 1526
$ws.Range("L94").Value = 1168.8462
$ws.Range("M94").Value = -1075
$ws.Range("N94").Value = -2070.8462
$ws.Range("H95").Value = 16507.25
$ws.Range("J95").Value = 16507.25
$ws.Range("L95").Value = 16507.25
$ws.Range("N95").Value = -21999.25
$ws.Range("H102").Value = 37466
$ws.Range("J102").Value = 38700
$ws.Range("L102").Value = 38700
$ws.Range("N102").Value = -43568
$ws.Range("H104").Value = 48263
$ws.Range("J104").Value = 48263
$ws.Range("L104").Value = 48263
$ws.Range("N104").Value = -53505
$ws.Range("H109").Value = 51055.332
$ws.Range("J109").Value = 51055.332
$ws.Range("L109").Value = 51055.332
$ws.Range("N109").Value = -53135.332
$ws.Range("H115").Value = 37800
$ws.Range("J115").Value = 37800
$ws.Range("L115").Value = 37800
$ws.Range("N115").Value = -40150
$ws.Range("H122").Value = 2154.95
$ws.Range("I122").Value = 1752.5714
$ws.Range("K122").Value = 5257.7142
$ws.Range("M122").Value = -2807.7142
$ws.Range("H132").Value = 3572.6086
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 7510.0625
$ws.Range("I134").Value = 7361.5
$ws.Range("J134").Value = 8550
$ws.Range("K134").Value = 22084.5
$ws.Range("L134").Value = 25650
$ws.Range("M134").Value = -19549.5
$ws.Range("N134").Value = -30720
$ws.Range("H136").Value = 9100.157999999999
$ws.Range("I136").Value = 3036.7273
$ws.Range("J136").Value = 17437.375
$ws.Range("K136").Value = 9110.1819
$ws.Range("L136").Value = 52312.125
$ws.Range("M136").Value = -6560.1819
$ws.Range("N136").Value = -57412.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.272728
$ws.Range("J2").Value = 54.125
$ws.Range("L2").Value = 324.75
$ws.Range("N2").Value = -550.75
$ws.Range("H37").Value = 117574.336
$ws.Range("J37").Value = 117574.336
$ws.Range("L37").Value = 352723.008
$ws.Range("N37").Value = -352947.008
$ws.Range("H113").Value = 2581.1
$ws.Range("I113").Value = 623
$ws.Range("J113").Value = 3070.625
$ws.Range("K113").Value = 1869
$ws.Range("L113").Value = 9211.875
$ws.Range("M113").Value = 301
$ws.Range("N113").Value = -13551.875
$ws.Range("H116").Value = 2769.2
$ws.Range("I116").Value = 2461.5
$ws.Range("K116").Value = 7384.5
$ws.Range("M116").Value = -3942.5
$ws.Range("H122").Value = 799
$ws.Range("J122").Value = 799
$ws.Range("L122").Value = 7191
$ws.Range("N122").Value = -12091
$ws.Range("H132").Value = 2616.3635
$ws.Range("I132").Value = 1729.5555
$ws.Range("J132").Value = 3230.3076
$ws.Range("K132").Value = 15565.9995
$ws.Range("L132").Value = 29072.7684
$ws.Range("M132").Value = -13035.9995
$ws.Range("N132").Value = -34132.7684
$ws.Range("H137").Value = 1740.4
$ws.Range("I137").Value = 1235.6666
$ws.Range("J137").Value = 2153.3635
$ws.Range("K137").Value = 3706.9998
$ws.Range("L137").Value = 6460.0905
$ws.Range("M137").Value = 1393.0002
$ws.Range("N137").Value = -16660.0905

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9001.25
$ws.Range("I70").Value = 8008
$ws.Range("K70").Value = 8008
$ws.Range("M70").Value = -7738
$ws.Range("H73").Value = 9001.25
$ws.Range("I73").Value = 8008
$ws.Range("K73").Value = 8008
$ws.Range("M73").Value = -7072
$ws.Range("H102").Value = 42116.785
$ws.Range("I102").Value = 56538.21
$ws.Range("K102").Value = 56538.21
$ws.Range("M102").Value = -54916.21
$ws.Range("H126").Value = 59647.875
$ws.Range("I126").Value = 67726.92999999999
$ws.Range("K126").Value = 203180.79
$ws.Range("M126").Value = -200710.79
$ws.Range("H132").Value = 3464.0454
$ws.Range("I132").Value = 4046.077
$ws.Range("J132").Value = 2623.3333
$ws.Range("K132").Value = 12138.231
$ws.Range("L132").Value = 7869.999899999999
$ws.Range("M132").Value = -9608.231
$ws.Range("N132").Value = -12929.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 959.12195
$ws.Range("I16").Value = 843.82855
$ws.Range("J16").Value = 1631.6666
$ws.Range("K16").Value = 843.82855
$ws.Range("L16").Value = 1631.6666
$ws.Range("M16").Value = -673.82855
$ws.Range("N16").Value = -1971.6666
$ws.Range("H40").Value = 3207.1428
$ws.Range("J40").Value = 4223.75
$ws.Range("L40").Value = 4223.75
$ws.Range("N40").Value = -4495.75
$ws.Range("H46").Value = 2579.9688
$ws.Range("I46").Value = 1041.3334
$ws.Range("J46").Value = 2935.0386
$ws.Range("K46").Value = 1041.3334
$ws.Range("L46").Value = 2935.0386
$ws.Range("M46").Value = -853.3334
$ws.Range("N46").Value = -3311.0386
$ws.Range("H59").Value = 29264.334
$ws.Range("J59").Value = 29264.334
$ws.Range("L59").Value = 29264.334
$ws.Range("N59").Value = -30572.334
$ws.Range("H61").Value = 18520270
$ws.Range("I61").Value = 31251506
$ws.Range("J61").Value = 2109.7273
$ws.Range("K61").Value = 31251506
$ws.Range("L61").Value = 2109.7273
$ws.Range("M61").Value = -31251304
$ws.Range("N61").Value = -2513.7273
$ws.Range("H68").Value = 4769.6
$ws.Range("I68").Value = 3712
$ws.Range("K68").Value = 3712
$ws.Range("M68").Value = -2963
$ws.Range("H71").Value = 4769.6
$ws.Range("I71").Value = 3712
$ws.Range("K71").Value = 18560
$ws.Range("M71").Value = -14816
$ws.Range("H100").Value = 2581.2778
$ws.Range("I100").Value = 1917.0769
$ws.Range("K100").Value = 1917.0769
$ws.Range("M100").Value = -1376.0769
$ws.Range("H113").Value = 18520270
$ws.Range("I113").Value = 31251506
$ws.Range("J113").Value = 2109.7273
$ws.Range("K113").Value = 31251506
$ws.Range("L113").Value = 2109.7273
$ws.Range("M113").Value = -31249336
$ws.Range("N113").Value = -6449.7273
$ws.Range("H122").Value = 3399.0244
$ws.Range("I122").Value = 3334.7297
$ws.Range("K122").Value = 10004.1891
$ws.Range("M122").Value = -7554.1891
$ws.Range("H132").Value = 2734.6572
$ws.Range("I132").Value = 2295.7585
$ws.Range("J132").Value = 4856
$ws.Range("K132").Value = 6887.2755
$ws.Range("L132").Value = 14568
$ws.Range("M132").Value = -4357.2755
$ws.Range("N132").Value = -19628
$ws.Range("H136").Value = 34103.965
$ws.Range("I136").Value = 2623.2104
$ws.Range("J136").Value = 93917.39999999999
$ws.Range("K136").Value = 7869.6312
$ws.Range("L136").Value = 281752.2
$ws.Range("M136").Value = -5319.6312
$ws.Range("N136").Value = -286852.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20911
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20911
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20911
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -21691
$ws.Range("H96").Value = 15670
$ws.Range("I96").Value = 15670
$ws.Range("K96").Value = 15670
$ws.Range("M96").Value = -14297
$ws.Range("H104").Value = 26807.857
$ws.Range("J104").Value = 26807.857
$ws.Range("L104").Value = 26807.857
$ws.Range("N104").Value = -33795.857
$ws.Range("H105").Value = 29666.334
$ws.Range("J105").Value = 29666.334
$ws.Range("L105").Value = 29666.334
$ws.Range("N105").Value = -36654.334
$ws.Range("H107").Value = 898.8
$ws.Range("I107").Value = 873.75
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 2621.25
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -701.25
$ws.Range("N107").Value = -6837
$ws.Range("H113").Value = 1790.7273
$ws.Range("I113").Value = 1769.8
$ws.Range("K113").Value = 5309.4
$ws.Range("M113").Value = -3139.4
$ws.Range("H126").Value = 52635560
$ws.Range("I126").Value = 71432424
$ws.Range("K126").Value = 214297272
$ws.Range("M126").Value = -214294802
$ws.Range("H132").Value = 3964.8235
$ws.Range("I132").Value = 4082.3076
$ws.Range("K132").Value = 12246.9228
$ws.Range("M132").Value = -9716.9228
$ws.Range("H133").Value = 65475
$ws.Range("J133").Value = 64546.75
$ws.Range("L133").Value = 64546.75
$ws.Range("N133").Value = -74666.75
$ws.Range("H136").Value = 46349.49
$ws.Range("I136").Value = 35183.668
$ws.Range("J136").Value = 67883.57000000001
$ws.Range("K136").Value = 105551.004
$ws.Range("L136").Value = 203650.71
$ws.Range("M136").Value = -103001.004
$ws.Range("N136").Value = -208750.71
